# Natmi following Dr Hou advice
# Update Jag1-Notch4 LR-pair stats: ligand/receptor-expressing cell counts
# changed from 1 to 3, with all dependent expression/specificity values
# recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.66920033333333
$ws.Range("H2").Value = 47.00760099999999
$ws.Range("I2").Value = 0.2925937299273087
$ws.Range("J2").Value = 0.2925937299273087
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 29.546424
$ws.Range("N2").Value = 88.639272
$ws.Range("O2").Value = 0.9033225104610835
$ws.Range("P2").Value = 0.9033225104610834
$ws.Range("Q2").Value = 462.968836789608
$ws.Range("R2").Value = 4166.719531106472
$ws.Range("S2").Value = 0.2643065026631087
$ws.Range("T2").Value = 0.2643065026631087

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.66920033333333
$ws.Range("H3").Value = 47.00760099999999
$ws.Range("I3").Value = 0.2925937299273087
$ws.Range("J3").Value = 0.2925937299273087
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.672785333333334
$ws.Range("N3").Value = 5.018356000000001
$ws.Range("O3").Value = 0.05114204841740398
$ws.Range("P3").Value = 0.05114204841740398
$ws.Range("Q3").Value = 26.21120850266178
$ws.Range("R3").Value = 235.900876523956
$ws.Range("S3").Value = 0.01496384270257125
$ws.Range("T3").Value = 0.01496384270257125

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.66920033333333
$ws.Range("H4").Value = 47.00760099999999
$ws.Range("I4").Value = 0.2925937299273087
$ws.Range("J4").Value = 0.2925937299273087
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.489401
$ws.Range("N4").Value = 4.468203
$ws.Range("O4").Value = 0.04553544112151264
$ws.Range("P4").Value = 0.04553544112151264
$ws.Range("Q4").Value = 23.337722645667
$ws.Range("R4").Value = 210.039503811003
$ws.Range("S4").Value = 0.01332338456162874
$ws.Range("T4").Value = 0.01332338456162874

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.47676966666667
$ws.Range("H5").Value = 49.43030900000001
$ws.Range("I5").Value = 0.3076736139282969
$ws.Range("J5").Value = 0.3076736139282968
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.546424
$ws.Range("N5").Value = 88.639272
$ws.Range("O5").Value = 0.9033225104610835
$ws.Range("P5").Value = 0.9033225104610834
$ws.Range("Q5").Value = 486.8296227216721
$ws.Range("R5").Value = 4381.466604495049
$ws.Range("S5").Value = 0.2779285013363433
$ws.Range("T5").Value = 0.2779285013363432

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.47676966666667
$ws.Range("H6").Value = 49.43030900000001
$ws.Range("I6").Value = 0.3076736139282969
$ws.Range("J6").Value = 0.3076736139282968
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.672785333333334
$ws.Range("N6").Value = 5.018356000000001
$ws.Range("O6").Value = 0.05114204841740398
$ws.Range("P6").Value = 0.05114204841740398
$ws.Range("Q6").Value = 27.56209863911156
$ws.Range("R6").Value = 248.0588877520041
$ws.Range("S6").Value = 0.01573505886027862
$ws.Range("T6").Value = 0.01573505886027862

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.47676966666667
$ws.Range("H7").Value = 49.43030900000001
$ws.Range("I7").Value = 0.3076736139282969
$ws.Range("J7").Value = 0.3076736139282968
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.489401
$ws.Range("N7").Value = 4.468203
$ws.Range("O7").Value = 0.04553544112151264
$ws.Range("P7").Value = 0.04553544112151264
$ws.Range("Q7").Value = 24.54051721830301
$ws.Range("R7").Value = 220.864654964727
$ws.Range("S7").Value = 0.01401005373167498
$ws.Range("T7").Value = 0.01401005373167497

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.406785
$ws.Range("H8").Value = 64.220355
$ws.Range("I8").Value = 0.3997326561443945
$ws.Range("J8").Value = 0.3997326561443944
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.546424
$ws.Range("N8").Value = 88.639272
$ws.Range("O8").Value = 0.9033225104610835
$ws.Range("P8").Value = 0.9033225104610834
$ws.Range("Q8").Value = 632.49394608684
$ws.Range("R8").Value = 5692.44551478156
$ws.Range("S8").Value = 0.3610875064616315
$ws.Range("T8").Value = 0.3610875064616313

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.406785
$ws.Range("H9").Value = 64.220355
$ws.Range("I9").Value = 0.3997326561443945
$ws.Range("J9").Value = 0.3997326561443944
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.672785333333334
$ws.Range("N9").Value = 5.018356000000001
$ws.Range("O9").Value = 0.05114204841740398
$ws.Range("P9").Value = 0.05114204841740398
$ws.Range("Q9").Value = 35.80895598182001
$ws.Range("R9").Value = 322.28060383638
$ws.Range("S9").Value = 0.02044314685455412
$ws.Range("T9").Value = 0.02044314685455412

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 21.406785
$ws.Range("H10").Value = 64.220355
$ws.Range("I10").Value = 0.3997326561443945
$ws.Range("J10").Value = 0.3997326561443944
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.489401
$ws.Range("N10").Value = 4.468203
$ws.Range("O10").Value = 0.04553544112151264
$ws.Range("P10").Value = 0.04553544112151264
$ws.Range("Q10").Value = 31.883286985785
$ws.Range("R10").Value = 286.949582872065
$ws.Range("S10").Value = 0.01820200282820893
$ws.Range("T10").Value = 0.01820200282820893
